$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-07 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-08 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("124×6=", $true, $false, $false, $false, $false, $true, 1, $false, "802×6=", 2) | Out-Null
$d.Content.Find.Execute("996×3=", $true, $false, $false, $false, $false, $true, 1, $false, "103×2=", 2) | Out-Null
$d.Content.Find.Execute("321×3=", $true, $false, $false, $false, $false, $true, 1, $false, "330×9=", 2) | Out-Null
$d.Content.Find.Execute("132×7=", $true, $false, $false, $false, $false, $true, 1, $false, "679×8=", 2) | Out-Null
$d.Content.Find.Execute("571×8=", $true, $false, $false, $false, $false, $true, 1, $false, "566×5=", 2) | Out-Null
$d.Content.Find.Execute("357×9=", $true, $false, $false, $false, $false, $true, 1, $false, "736×6=", 2) | Out-Null
$d.Content.Find.Execute("386×4=", $true, $false, $false, $false, $false, $true, 1, $false, "622×5=", 2) | Out-Null
$d.Content.Find.Execute("759×6=", $true, $false, $false, $false, $false, $true, 1, $false, "560×4=", 2) | Out-Null
$d.Content.Find.Execute("624×6=", $true, $false, $false, $false, $false, $true, 1, $false, "730×5=", 2) | Out-Null
$d.Content.Find.Execute("744×2=", $true, $false, $false, $false, $false, $true, 1, $false, "716×9=", 2) | Out-Null
$d.Content.Find.Execute("409×3=", $true, $false, $false, $false, $false, $true, 1, $false, "140×8=", 2) | Out-Null
$d.Content.Find.Execute("918×6=", $true, $false, $false, $false, $false, $true, 1, $false, "827×2=", 2) | Out-Null
$d.Content.Find.Execute("675×9=", $true, $false, $false, $false, $false, $true, 1, $false, "641×4=", 2) | Out-Null
$d.Content.Find.Execute("356×7=", $true, $false, $false, $false, $false, $true, 1, $false, "463×7=", 2) | Out-Null
$d.Content.Find.Execute("739×6=", $true, $false, $false, $false, $false, $true, 1, $false, "981×6=", 2) | Out-Null
$d.Content.Find.Execute("631×4=", $true, $false, $false, $false, $false, $true, 1, $false, "900×6=", 2) | Out-Null
$d.Content.Find.Execute("452×9=", $true, $false, $false, $false, $false, $true, 1, $false, "407×9=", 2) | Out-Null
$d.Content.Find.Execute("402×5=", $true, $false, $false, $false, $false, $true, 1, $false, "955×6=", 2) | Out-Null
$d.Content.Find.Execute("681×8=", $true, $false, $false, $false, $false, $true, 1, $false, "581×6=", 2) | Out-Null
$d.Content.Find.Execute("417×3=", $true, $false, $false, $false, $false, $true, 1, $false, "359×4=", 2) | Out-Null
$d.Content.Find.Execute("823×6=", $true, $false, $false, $false, $false, $true, 1, $false, "824×2=", 2) | Out-Null
$d.Content.Find.Execute("317×4=", $true, $false, $false, $false, $false, $true, 1, $false, "931×5=", 2) | Out-Null
$d.Content.Find.Execute("242×3=", $true, $false, $false, $false, $false, $true, 1, $false, "298×5=", 2) | Out-Null
$d.Content.Find.Execute("861×9=", $true, $false, $false, $false, $false, $true, 1, $false, "214×4=", 2) | Out-Null
$d.Content.Find.Execute("596×8=", $true, $false, $false, $false, $false, $true, 1, $false, "892×9=", 2) | Out-Null
